$d = $word.ActiveDocument

[void]$d.Content.Find.Execute('71×76=5396', $true, $false, $false, $false, $false, $true, 1, $false, '71×70=4970', 2)
[void]$d.Content.Find.Execute('67×15=1005', $true, $false, $false, $false, $false, $true, 1, $false, '37×98=3626', 2)
[void]$d.Content.Find.Execute('74×75=5550', $true, $false, $false, $false, $false, $true, 1, $false, '30×79=2370', 2)
[void]$d.Content.Find.Execute('94×41=3854', $true, $false, $false, $false, $false, $true, 1, $false, '12×17=204', 2)
[void]$d.Content.Find.Execute('42×73=3066', $true, $false, $false, $false, $false, $true, 1, $false, '75×43=3225', 2)
[void]$d.Content.Find.Execute('76×87=6612', $true, $false, $false, $false, $false, $true, 1, $false, '95×91=8645', 2)
[void]$d.Content.Find.Execute('78×79=6162', $true, $false, $false, $false, $false, $true, 1, $false, '72×47=3384', 2)
[void]$d.Content.Find.Execute('68×78=5304', $true, $false, $false, $false, $false, $true, 1, $false, '49×78=3822', 2)
[void]$d.Content.Find.Execute('40×66=2640', $true, $false, $false, $false, $false, $true, 1, $false, '81×45=3645', 2)
[void]$d.Content.Find.Execute('32×76=2432', $true, $false, $false, $false, $false, $true, 1, $false, '99×81=8019', 2)
[void]$d.Content.Find.Execute('39×77=3003', $true, $false, $false, $false, $false, $true, 1, $false, '71×85=6035', 2)
[void]$d.Content.Find.Execute('98×11=1078', $true, $false, $false, $false, $false, $true, 1, $false, '56×12=672', 2)
[void]$d.Content.Find.Execute('62×14=868', $true, $false, $false, $false, $false, $true, 1, $false, '20×13=260', 2)
[void]$d.Content.Find.Execute('22×66=1452', $true, $false, $false, $false, $false, $true, 1, $false, '29×83=2407', 2)
[void]$d.Content.Find.Execute('14×45=630', $true, $false, $false, $false, $false, $true, 1, $false, '84×94=7896', 2)
[void]$d.Content.Find.Execute('70×44=3080', $true, $false, $false, $false, $false, $true, 1, $false, '59×31=1829', 2)
[void]$d.Content.Find.Execute('61×97=5917', $true, $false, $false, $false, $false, $true, 1, $false, '21×34=714', 2)
[void]$d.Content.Find.Execute('19×61=1159', $true, $false, $false, $false, $false, $true, 1, $false, '50×93=4650', 2)
[void]$d.Content.Find.Execute('78×64=4992', $true, $false, $false, $false, $false, $true, 1, $false, '21×23=483', 2)
[void]$d.Content.Find.Execute('15×28=420', $true, $false, $false, $false, $false, $true, 1, $false, '14×40=560', 2)
[void]$d.Content.Find.Execute('94×23=2162', $true, $false, $false, $false, $false, $true, 1, $false, '67×85=5695', 2)
[void]$d.Content.Find.Execute('48×14=672', $true, $false, $false, $false, $false, $true, 1, $false, '74×52=3848', 2)
[void]$d.Content.Find.Execute('36×10=360', $true, $false, $false, $false, $false, $true, 1, $false, '63×72=4536', 2)
[void]$d.Content.Find.Execute('60×55=3300', $true, $false, $false, $false, $false, $true, 1, $false, '45×14=630', 2)
[void]$d.Content.Find.Execute('38×88=3344', $true, $false, $false, $false, $false, $true, 1, $false, '52×20=1040', 2)
[void]$d.Content.Find.Execute('56×28=1568', $true, $false, $false, $false, $false, $true, 1, $false, '52×44=2288', 2)
[void]$d.Content.Find.Execute('80×29=2320', $true, $false, $false, $false, $false, $true, 1, $false, '58×33=1914', 2)
[void]$d.Content.Find.Execute('95×65=6175', $true, $false, $false, $false, $false, $true, 1, $false, '88×86=7568', 2)
[void]$d.Content.Find.Execute('86×21=1806', $true, $false, $false, $false, $false, $true, 1, $false, '52×22=1144', 2)
[void]$d.Content.Find.Execute('48×96=4608', $true, $false, $false, $false, $false, $true, 1, $false, '39×89=3471', 2)
[void]$d.Content.Find.Execute('34×53=1802', $true, $false, $false, $false, $false, $true, 1, $false, '15×66=990', 2)
[void]$d.Content.Find.Execute('44×83=3652', $true, $false, $false, $false, $false, $true, 1, $false, '68×88=5984', 2)
[void]$d.Content.Find.Execute('87×86=7482', $true, $false, $false, $false, $false, $true, 1, $false, '21×97=2037', 2)
[void]$d.Content.Find.Execute('21×89=1869', $true, $false, $false, $false, $false, $true, 1, $false, '99×32=3168', 2)
[void]$d.Content.Find.Execute('40×77=3080', $true, $false, $false, $false, $false, $true, 1, $false, '42×63=2646', 2)
[void]$d.Content.Find.Execute('98×62=6076', $true, $false, $false, $false, $false, $true, 1, $false, '97×66=6402', 2)
[void]$d.Content.Find.Execute('36×60=2160', $true, $false, $false, $false, $false, $true, 1, $false, '96×79=7584', 2)
[void]$d.Content.Find.Execute('90×45=4050', $true, $false, $false, $false, $false, $true, 1, $false, '13×83=1079', 2)
[void]$d.Content.Find.Execute('83×73=6059', $true, $false, $false, $false, $false, $true, 1, $false, '72×60=4320', 2)
[void]$d.Content.Find.Execute('87×24=2088', $true, $false, $false, $false, $false, $true, 1, $false, '50×71=3550', 2)
[void]$d.Content.Find.Execute('42×95=3990', $true, $false, $false, $false, $false, $true, 1, $false, '17×10=170', 2)
[void]$d.Content.Find.Execute('61×36=2196', $true, $false, $false, $false, $false, $true, 1, $false, '61×87=5307', 2)
[void]$d.Content.Find.Execute('74×45=3330', $true, $false, $false, $false, $false, $true, 1, $false, '43×48=2064', 2)
[void]$d.Content.Find.Execute('100×30=3000', $true, $false, $false, $false, $false, $true, 1, $false, '80×23=1840', 2)
[void]$d.Content.Find.Execute('31×26=806', $true, $false, $false, $false, $false, $true, 1, $false, '22×57=1254', 2)
[void]$d.Content.Find.Execute('24×70=1680', $true, $false, $false, $false, $false, $true, 1, $false, '18×33=594', 2)
[void]$d.Content.Find.Execute('74×51=3774', $true, $false, $false, $false, $false, $true, 1, $false, '80×89=7120', 2)
[void]$d.Content.Find.Execute('47×44=2068', $true, $false, $false, $false, $false, $true, 1, $false, '13×98=1274', 2)
[void]$d.Content.Find.Execute('46×27=1242', $true, $false, $false, $false, $false, $true, 1, $false, '10×95=950', 2)
[void]$d.Content.Find.Execute('68×27=1836', $true, $false, $false, $false, $false, $true, 1, $false, '22×14=308', 2)
[void]$d.Content.Find.Execute('27×70=1890', $true, $false, $false, $false, $false, $true, 1, $false, '91×79=7189', 2)
[void]$d.Content.Find.Execute('72×77=5544', $true, $false, $false, $false, $false, $true, 1, $false, '85×38=3230', 2)
[void]$d.Content.Find.Execute('53×32=1696', $true, $false, $false, $false, $false, $true, 1, $false, '95×10=950', 2)
[void]$d.Content.Find.Execute('94×70=6580', $true, $false, $false, $false, $false, $true, 1, $false, '11×69=759', 2)
[void]$d.Content.Find.Execute('21×48=1008', $true, $false, $false, $false, $false, $true, 1, $false, '41×19=779', 2)
[void]$d.Content.Find.Execute('87×43=3741', $true, $false, $false, $false, $false, $true, 1, $false, '79×54=4266', 2)
[void]$d.Content.Find.Execute('32×83=2656', $true, $false, $false, $false, $false, $true, 1, $false, '86×81=6966', 2)
[void]$d.Content.Find.Execute('95×14=1330', $true, $false, $false, $false, $false, $true, 1, $false, '94×94=8836', 2)
[void]$d.Content.Find.Execute('71×38=2698', $true, $false, $false, $false, $false, $true, 1, $false, '69×93=6417', 2)
[void]$d.Content.Find.Execute('79×98=7742', $true, $false, $false, $false, $false, $true, 1, $false, '30×89=2670', 2)
[void]$d.Content.Find.Execute('69×41=2829', $true, $false, $false, $false, $false, $true, 1, $false, '76×79=6004', 2)
[void]$d.Content.Find.Execute('24×40=960', $true, $false, $false, $false, $false, $true, 1, $false, '74×65=4810', 2)
[void]$d.Content.Find.Execute('30×61=1830', $true, $false, $false, $false, $false, $true, 1, $false, '61×87=5307', 2)
[void]$d.Content.Find.Execute('30×73=2190', $true, $false, $false, $false, $false, $true, 1, $false, '13×77=1001', 2)
[void]$d.Content.Find.Execute('13×69=897', $true, $false, $false, $false, $false, $true, 1, $false, '24×21=504', 2)
[void]$d.Content.Find.Execute('15×21=315', $true, $false, $false, $false, $false, $true, 1, $false, '73×57=4161', 2)
[void]$d.Content.Find.Execute('32×68=2176', $true, $false, $false, $false, $false, $true, 1, $false, '74×30=2220', 2)
[void]$d.Content.Find.Execute('90×62=5580', $true, $false, $false, $false, $false, $true, 1, $false, '24×35=840', 2)
[void]$d.Content.Find.Execute('62×47=2914', $true, $false, $false, $false, $false, $true, 1, $false, '65×45=2925', 2)
[void]$d.Content.Find.Execute('54×26=1404', $true, $false, $false, $false, $false, $true, 1, $false, '30×96=2880', 2)
[void]$d.Content.Find.Execute('42×50=2100', $true, $false, $false, $false, $false, $true, 1, $false, '72×80=5760', 2)
[void]$d.Content.Find.Execute('42×29=1218', $true, $false, $false, $false, $false, $true, 1, $false, '63×40=2520', 2)
[void]$d.Content.Find.Execute('12×34=408', $true, $false, $false, $false, $false, $true, 1, $false, '28×62=1736', 2)
[void]$d.Content.Find.Execute('22×34=748', $true, $false, $false, $false, $false, $true, 1, $false, '38×99=3762', 2)
[void]$d.Content.Find.Execute('41×60=2460', $true, $false, $false, $false, $false, $true, 1, $false, '93×18=1674', 2)
[void]$d.Content.Find.Execute('92×82=7544', $true, $false, $false, $false, $false, $true, 1, $false, '81×38=3078', 2)
[void]$d.Content.Find.Execute('76×76=5776', $true, $false, $false, $false, $false, $true, 1, $false, '47×26=1222', 2)
[void]$d.Content.Find.Execute('70×32=2240', $true, $false, $false, $false, $false, $true, 1, $false, '19×70=1330', 2)
[void]$d.Content.Find.Execute('16×53=848', $true, $false, $false, $false, $false, $true, 1, $false, '59×87=5133', 2)
[void]$d.Content.Find.Execute('98×57=5586', $true, $false, $false, $false, $false, $true, 1, $false, '95×90=8550', 2)
[void]$d.Content.Find.Execute('74×38=2812', $true, $false, $false, $false, $false, $true, 1, $false, '42×18=756', 2)
[void]$d.Content.Find.Execute('56×25=1400', $true, $false, $false, $false, $false, $true, 1, $false, '30×31=930', 2)
[void]$d.Content.Find.Execute('54×100=5400', $true, $false, $false, $false, $false, $true, 1, $false, '23×65=1495', 2)
[void]$d.Content.Find.Execute('52×27=1404', $true, $false, $false, $false, $false, $true, 1, $false, '55×64=3520', 2)
[void]$d.Content.Find.Execute('44×22=968', $true, $false, $false, $false, $false, $true, 1, $false, '72×86=6192', 2)
[void]$d.Content.Find.Execute('87×33=2871', $true, $false, $false, $false, $false, $true, 1, $false, '81×26=2106', 2)
[void]$d.Content.Find.Execute('11×13=143', $true, $false, $false, $false, $false, $true, 1, $false, '59×58=3422', 2)
[void]$d.Content.Find.Execute('89×78=6942', $true, $false, $false, $false, $false, $true, 1, $false, '81×39=3159', 2)
[void]$d.Content.Find.Execute('38×82=3116', $true, $false, $false, $false, $false, $true, 1, $false, '53×48=2544', 2)
[void]$d.Content.Find.Execute('50×60=3000', $true, $false, $false, $false, $false, $true, 1, $false, '66×56=3696', 2)
[void]$d.Content.Find.Execute('42×83=3486', $true, $false, $false, $false, $false, $true, 1, $false, '23×63=1449', 2)
[void]$d.Content.Find.Execute('11×45=495', $true, $false, $false, $false, $false, $true, 1, $false, '97×87=8439', 2)
[void]$d.Content.Find.Execute('67×69=4623', $true, $false, $false, $false, $false, $true, 1, $false, '24×43=1032', 2)
[void]$d.Content.Find.Execute('12×57=684', $true, $false, $false, $false, $false, $true, 1, $false, '66×21=1386', 2)
[void]$d.Content.Find.Execute('22×72=1584', $true, $false, $false, $false, $false, $true, 1, $false, '40×90=3600', 2)
[void]$d.Content.Find.Execute('44×92=4048', $true, $false, $false, $false, $false, $true, 1, $false, '21×33=693', 2)
[void]$d.Content.Find.Execute('11×43=473', $true, $false, $false, $false, $false, $true, 1, $false, '63×51=3213', 2)
[void]$d.Content.Find.Execute('55×22=1210', $true, $false, $false, $false, $false, $true, 1, $false, '52×67=3484', 2)
[void]$d.Content.Find.Execute('83×68=5644', $true, $false, $false, $false, $false, $true, 1, $false, '21×36=756', 2)
[void]$d.Content.Find.Execute('30×12=360', $true, $false, $false, $false, $false, $true, 1, $false, '24×14=336', 2)
